$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 371.04544  # H33: was 393.35
$ws.Cells.Item(33, 9).Value = 298.2  # I33: was 314.8889
$ws.Cells.Item(33, 11).Value = 298.2  # K33: was 314.8889
$ws.Cells.Item(33, 13).Value = -69.19999999999999  # M33: was -85.88889999999998
$ws.Cells.Item(103, 8).Value = 527  # H103: was 650.5714
$ws.Cells.Item(103, 10).Value = 480.5  # J103: was 609.1667
$ws.Cells.Item(103, 12).Value = 1441.5  # L103: was 1827.5001
$ws.Cells.Item(103, 14).Value = -2613.5  # N103: was -2999.5001
$ws.Cells.Item(112, 8).Value = 3788.4324  # H112: was 3838.139
$ws.Cells.Item(112, 9).Value = 6674.4  # I112: was 6714.4
$ws.Cells.Item(112, 10).Value = 3337.5  # J112: was 3374.2258
$ws.Cells.Item(112, 11).Value = 20023.2  # K112: was 20143.2
$ws.Cells.Item(112, 12).Value = 10012.5  # L112: was 10122.6774
$ws.Cells.Item(112, 13).Value = -18915.2  # M112: was -19035.2
$ws.Cells.Item(112, 14).Value = -12228.5  # N112: was -12338.6774
$ws.Cells.Item(127, 8).Value = 4397.16  # H127: was 4937.864
$ws.Cells.Item(127, 9).Value = 4397.16  # I127: was 4937.864
$ws.Cells.Item(127, 11).Value = 13191.48  # K127: was 14813.592
$ws.Cells.Item(127, 13).Value = -8231.48  # M127: was -9853.591999999999
$ws.Cells.Item(132, 8).Value = 503345.62  # H132: was 544924.6
$ws.Cells.Item(132, 9).Value = 683120.4399999999  # I132: was 720860.75
$ws.Cells.Item(132, 10).Value = 15385.429  # J132: was 17116.334
$ws.Cells.Item(132, 11).Value = 2049361.32  # K132: was 2162582.25
$ws.Cells.Item(132, 12).Value = 46156.287  # L132: was 51349.00199999999
$ws.Cells.Item(132, 13).Value = -2046831.32  # M132: was -2160052.25
$ws.Cells.Item(132, 14).Value = -51216.287  # N132: was -56409.00199999999
$ws.Cells.Item(137, 8).Value = 7502.32  # H137: was 7518.48
$ws.Cells.Item(137, 9).Value = 4076.8125  # I137: was 4255.6
$ws.Cells.Item(137, 10).Value = 13592.111  # J137: was 12412.8
$ws.Cells.Item(137, 11).Value = 12230.4375  # K137: was 12766.8
$ws.Cells.Item(137, 12).Value = 40776.333  # L137: was 37238.39999999999
$ws.Cells.Item(137, 13).Value = -9680.4375  # M137: was -10216.8
$ws.Cells.Item(137, 14).Value = -45876.333  # N137: was -42338.39999999999
$ws.Cells.Item(138, 8).Value = 3870.1162  # H138: was 3882.0234
$ws.Cells.Item(138, 10).Value = 4156.507  # J138: was 4174.5415
$ws.Cells.Item(138, 12).Value = 12469.521  # L138: was 12523.6245
$ws.Cells.Item(138, 14).Value = -22749.521  # N138: was -22803.6245

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2224479  # H32: was 2366447.8
$ws.Cells.Item(32, 9).Value = 1134.4634  # I32: was 1200.421
$ws.Cells.Item(32, 11).Value = 1134.4634  # K32: was 1200.421
$ws.Cells.Item(32, 13).Value = -847.4634000000001  # M32: was -913.421
$ws.Cells.Item(45, 8).Value = 2351.0625  # H45: was 2442.4666
$ws.Cells.Item(45, 9).Value = 1955.2307  # I45: was 2036.5
$ws.Cells.Item(45, 11).Value = 1955.2307  # K45: was 2036.5
$ws.Cells.Item(45, 13).Value = -1578.2307  # M45: was -1659.5
$ws.Cells.Item(122, 8).Value = 4877.7354  # H122: was 4919.485
$ws.Cells.Item(122, 10).Value = 5915.5713  # J122: was 6101.385
$ws.Cells.Item(122, 12).Value = 17746.7139  # L122: was 18304.155
$ws.Cells.Item(122, 14).Value = -22646.7139  # N122: was -23204.155
$ws.Cells.Item(132, 8).Value = 559283.1  # H132: was 592768.2
$ws.Cells.Item(132, 9).Value = 610596.2  # I132: was 652146
$ws.Cells.Item(132, 11).Value = 1831788.6  # K132: was 1956438
$ws.Cells.Item(132, 13).Value = -1829258.6  # M132: was -1953908

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 5973215  # H107: was 6769274
$ws.Cells.Item(107, 9).Value = 6769350.5  # I107: was 7810362.5
$ws.Cells.Item(107, 11).Value = 6769350.5  # K107: was 7810362.5
$ws.Cells.Item(107, 13).Value = -6767430.5  # M107: was -7808442.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(29, 8).Value = 11999.5  # H29: was 0
$ws.Cells.Item(29, 10).Value = 11999.5  # J29: was 0
$ws.Cells.Item(29, 12).Value = 11999.5  # L29: was 0
$ws.Cells.Item(29, 14).Value = -12585.5  # N29: was None
$ws.Cells.Item(58, 8).Value = 55565150  # H58: was 58827910
$ws.Cells.Item(58, 9).Value = 71434616  # I58: was 76925740
$ws.Cells.Item(58, 10).Value = 22020  # J58: was 9965
$ws.Cells.Item(58, 11).Value = 71434616  # K58: was 76925740
$ws.Cells.Item(58, 12).Value = 22020  # L58: was 9965
$ws.Cells.Item(58, 13).Value = -71434413  # M58: was -76925537
$ws.Cells.Item(58, 14).Value = -22426  # N58: was -10371
$ws.Cells.Item(97, 8).Value = 0  # H97: was 41374.75
$ws.Cells.Item(97, 9).Value = 0  # I97: was 32999.5
$ws.Cells.Item(97, 10).Value = 0  # J97: was 49750
$ws.Cells.Item(97, 11).Value = 0  # K97: was 32999.5
$ws.Cells.Item(97, 12).Value = 0  # L97: was 49750
$ws.Cells.Item(97, 13).ClearContents()  # M97: was -32008.5
$ws.Cells.Item(97, 14).ClearContents()  # N97: was -51732
$ws.Cells.Item(105, 8).Value = 30306964  # H105: was 30306966
$ws.Cells.Item(105, 10).Value = 5546.5  # J105: was 5562.1665
$ws.Cells.Item(105, 12).Value = 5546.5  # L105: was 5562.1665
$ws.Cells.Item(105, 14).Value = -9040.5  # N105: was -9056.166499999999
$ws.Cells.Item(136, 8).Value = 55565150  # H136: was 58827910
$ws.Cells.Item(136, 9).Value = 71434616  # I136: was 76925740
$ws.Cells.Item(136, 10).Value = 22020  # J136: was 9965
$ws.Cells.Item(136, 11).Value = 214303848  # K136: was 230777220
$ws.Cells.Item(136, 12).Value = 66060  # L136: was 29895
$ws.Cells.Item(136, 13).Value = -214301298  # M136: was -230774670
$ws.Cells.Item(136, 14).Value = -71160  # N136: was -34995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 175.45454  # H2: was 216.12
$ws.Cells.Item(2, 9).Value = 149.66667  # I2: was 200.5
$ws.Cells.Item(2, 10).Value = 196.94444  # J2: was 226.53334
$ws.Cells.Item(2, 11).Value = 898.0000200000001  # K2: was 1203
$ws.Cells.Item(2, 12).Value = 1181.66664  # L2: was 1359.20004
$ws.Cells.Item(2, 13).Value = -785.0000200000001  # M2: was -1090
$ws.Cells.Item(2, 14).Value = -1407.66664  # N2: was -1585.20004
$ws.Cells.Item(34, 8).Value = 4448.467  # H34: was 4232.9375
$ws.Cells.Item(34, 10).Value = 9438.286  # J34: was 8383.5
$ws.Cells.Item(34, 12).Value = 28314.858  # L34: was 25150.5
$ws.Cells.Item(34, 14).Value = -28482.858  # N34: was -25318.5
$ws.Cells.Item(55, 8).Value = 18088.375  # H55: was 18962.375
$ws.Cells.Item(55, 9).Value = 4502.6665  # I55: was 5500
$ws.Cells.Item(55, 10).Value = 26239.8  # J55: was 23449.834
$ws.Cells.Item(55, 11).Value = 13507.9995  # K55: was 16500
$ws.Cells.Item(55, 12).Value = 78719.39999999999  # L55: was 70349.50199999999
$ws.Cells.Item(55, 13).Value = -13330.9995  # M55: was -16323
$ws.Cells.Item(55, 14).Value = -79073.39999999999  # N55: was -70703.50199999999
$ws.Cells.Item(62, 8).Value = 15166  # H62: was 15166.223
$ws.Cells.Item(62, 10).Value = 18083  # J62: was 18083.334
$ws.Cells.Item(62, 12).Value = 54249  # L62: was 54250.00199999999
$ws.Cells.Item(62, 14).Value = -55621  # N62: was -55622.00199999999
$ws.Cells.Item(65, 8).Value = 15166  # H65: was 15166.223
$ws.Cells.Item(65, 10).Value = 18083  # J65: was 18083.334
$ws.Cells.Item(65, 12).Value = 162747  # L65: was 162750.006
$ws.Cells.Item(65, 14).Value = -169611  # N65: was -169614.006
$ws.Cells.Item(92, 8).Value = 226  # H92: was 190.25
$ws.Cells.Item(92, 9).Value = 139.6  # I92: was 154
$ws.Cells.Item(92, 10).Value = 442  # J92: was 444
$ws.Cells.Item(92, 11).Value = 418.8  # K92: was 462
$ws.Cells.Item(92, 12).Value = 1326  # L92: was 1332
$ws.Cells.Item(92, 13).Value = 829.2  # M92: was 786
$ws.Cells.Item(92, 14).Value = -3822  # N92: was -3828
$ws.Cells.Item(107, 8).Value = 4274.615  # H107: was 4248.9487
$ws.Cells.Item(107, 10).Value = 6198.56  # J107: was 6158.52
$ws.Cells.Item(107, 12).Value = 18595.68  # L107: was 18475.56
$ws.Cells.Item(107, 14).Value = -22435.68  # N107: was -22315.56
$ws.Cells.Item(114, 8).Value = 636.8333  # H114: was 908.6667
$ws.Cells.Item(114, 9).Value = 636.8333  # I114: was 908.6667
$ws.Cells.Item(114, 11).Value = 1910.4999  # K114: was 2726.0001
$ws.Cells.Item(114, 13).Value = 1343.5001  # M114: was 527.9998999999998
$ws.Cells.Item(122, 8).Value = 119243.47  # H122: was 119225.47
$ws.Cells.Item(122, 9).Value = 415  # I122: was 398.75
$ws.Cells.Item(122, 10).Value = 160090.75  # J122: was 165222.9
$ws.Cells.Item(122, 11).Value = 3735  # K122: was 3588.75
$ws.Cells.Item(122, 12).Value = 1440816.75  # L122: was 1487006.1
$ws.Cells.Item(122, 13).Value = -1285  # M122: was -1138.75
$ws.Cells.Item(122, 14).Value = -1445716.75  # N122: was -1491906.1
$ws.Cells.Item(124, 8).Value = 11371  # H124: was 11496.4
$ws.Cells.Item(124, 10).Value = 20998.5  # J124: was 17998.334
$ws.Cells.Item(124, 12).Value = 62995.5  # L124: was 53995.00199999999
$ws.Cells.Item(124, 14).Value = -72815.5  # N124: was -63815.00199999999
$ws.Cells.Item(131, 8).Value = 54173380  # H131: was 54173384
$ws.Cells.Item(131, 10).Value = 47633776  # J131: was 47633790
$ws.Cells.Item(131, 12).Value = 142901328  # L131: was 142901370
$ws.Cells.Item(131, 14).Value = -142911408  # N131: was -142911450
$ws.Cells.Item(132, 8).Value = 2646.25  # H132: was 2726.3872
$ws.Cells.Item(132, 10).Value = 2765.2068  # J132: was 2858.1785
$ws.Cells.Item(132, 12).Value = 24886.8612  # L132: was 25723.6065
$ws.Cells.Item(132, 14).Value = -29946.8612  # N132: was -30783.6065

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 2299.5  # H7: was 6150
$ws.Cells.Item(7, 9).Value = 0  # I7: was 10000
$ws.Cells.Item(7, 10).Value = 2299.5  # J7: was 2300
$ws.Cells.Item(7, 11).Value = 0  # K7: was 10000
$ws.Cells.Item(7, 12).Value = 2299.5  # L7: was 2300
$ws.Cells.Item(7, 13).ClearContents()  # M7: was -9888
$ws.Cells.Item(7, 14).Value = -2523.5  # N7: was -2524
$ws.Cells.Item(8, 8).Value = 2299.5  # H8: was 6150
$ws.Cells.Item(8, 9).Value = 0  # I8: was 10000
$ws.Cells.Item(8, 10).Value = 2299.5  # J8: was 2300
$ws.Cells.Item(8, 11).Value = 0  # K8: was 10000
$ws.Cells.Item(8, 12).Value = 2299.5  # L8: was 2300
$ws.Cells.Item(8, 13).ClearContents()  # M8: was -9861
$ws.Cells.Item(8, 14).Value = -2577.5  # N8: was -2578
$ws.Cells.Item(11, 8).Value = 13339081  # H11: was 17015898
$ws.Cells.Item(11, 9).Value = 7785442.5  # I11: was 11451285
$ws.Cells.Item(11, 11).Value = 7785442.5  # K11: was 11451285
$ws.Cells.Item(11, 13).Value = -7785303.5  # M11: was -11451146
$ws.Cells.Item(113, 8).Value = 5949.8237  # H113: was 5708.6113
$ws.Cells.Item(113, 9).Value = 2047.4445  # I113: was 2003.5
$ws.Cells.Item(113, 11).Value = 2047.4445  # K113: was 2003.5
$ws.Cells.Item(113, 13).Value = 122.5554999999999  # M113: was 166.5
$ws.Cells.Item(122, 9).Value = 4563.65  # I122: was 4751.2104
$ws.Cells.Item(122, 10).Value = 15375  # J122: was 10583.333
$ws.Cells.Item(122, 11).Value = 13690.95  # K122: was 14253.6312
$ws.Cells.Item(122, 12).Value = 46125  # L122: was 31749.999
$ws.Cells.Item(122, 13).Value = -11240.95  # M122: was -11803.6312
$ws.Cells.Item(122, 14).Value = -51025  # N122: was -36649.999
$ws.Cells.Item(126, 8).Value = 50008804  # H126: was 55564970
$ws.Cells.Item(126, 9).Value = 125003080  # I126: was 166669660
$ws.Cells.Item(126, 11).Value = 375009240  # K126: was 500008980
$ws.Cells.Item(126, 13).Value = -375006770  # M126: was -500006510
$ws.Cells.Item(132, 8).Value = 58847950  # H132: was 58847930
$ws.Cells.Item(132, 9).Value = 100040504  # I132: was 90945970
$ws.Cells.Item(132, 10).Value = 1446.5714  # J132: was 1521
$ws.Cells.Item(132, 11).Value = 300121512  # K132: was 272837910
$ws.Cells.Item(132, 12).Value = 4339.7142  # L132: was 4563
$ws.Cells.Item(132, 13).Value = -300118982  # M132: was -272835380
$ws.Cells.Item(132, 14).Value = -9399.7142  # N132: was -9623

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 13946.8  # H7: was 12747.429
$ws.Cells.Item(7, 9).Value = 0  # I7: was 9749
$ws.Cells.Item(7, 11).Value = 0  # K7: was 9749
$ws.Cells.Item(7, 13).ClearContents()  # M7: was -9637
$ws.Cells.Item(40, 8).Value = 6161.2085  # H40: was 5175.4707
$ws.Cells.Item(40, 9).Value = 4493.125  # I40: was 4599.067
$ws.Cells.Item(40, 10).Value = 9497.375  # J40: was 9498.5
$ws.Cells.Item(40, 11).Value = 4493.125  # K40: was 4599.067
$ws.Cells.Item(40, 12).Value = 9497.375  # L40: was 9498.5
$ws.Cells.Item(40, 13).Value = -4357.125  # M40: was -4463.067
$ws.Cells.Item(40, 14).Value = -9769.375  # N40: was -9770.5
$ws.Cells.Item(46, 8).Value = 55560870  # H46: was 55557550
$ws.Cells.Item(46, 9).Value = 999  # I46: was 1033
$ws.Cells.Item(46, 10).Value = 62505856  # J46: was 83335816
$ws.Cells.Item(46, 11).Value = 999  # K46: was 1033
$ws.Cells.Item(46, 12).Value = 62505856  # L46: was 83335816
$ws.Cells.Item(46, 13).Value = -811  # M46: was -845
$ws.Cells.Item(46, 14).Value = -62506232  # N46: was -83336192
$ws.Cells.Item(55, 8).Value = 2372.6  # H55: was 2288.9614
$ws.Cells.Item(55, 9).Value = 718.0714  # I55: was 683.4
$ws.Cells.Item(55, 11).Value = 718.0714  # K55: was 683.4
$ws.Cells.Item(55, 13).Value = -545.0714  # M55: was -510.4
$ws.Cells.Item(68, 8).Value = 1266.5  # H68: was 1424.75
$ws.Cells.Item(68, 9).Value = 1199.6666  # I68: was 1299.6666
$ws.Cells.Item(68, 10).Value = 1333.3334  # J68: was 1800
$ws.Cells.Item(68, 11).Value = 1199.6666  # K68: was 1299.6666
$ws.Cells.Item(68, 12).Value = 1333.3334  # L68: was 1800
$ws.Cells.Item(68, 13).Value = -450.6666  # M68: was -550.6666
$ws.Cells.Item(68, 14).Value = -2831.3334  # N68: was -3298
$ws.Cells.Item(71, 8).Value = 1266.5  # H71: was 1424.75
$ws.Cells.Item(71, 9).Value = 1199.6666  # I71: was 1299.6666
$ws.Cells.Item(71, 10).Value = 1333.3334  # J71: was 1800
$ws.Cells.Item(71, 11).Value = 5998.333000000001  # K71: was 6498.333000000001
$ws.Cells.Item(71, 12).Value = 6666.666999999999  # L71: was 9000
$ws.Cells.Item(71, 13).Value = -2254.333000000001  # M71: was -2754.333000000001
$ws.Cells.Item(71, 14).Value = -14154.667  # N71: was -16488
$ws.Cells.Item(93, 8).Value = 1205.9  # H93: was 1192.7097
$ws.Cells.Item(93, 9).Value = 1191  # I93: was 1169.1111
$ws.Cells.Item(93, 11).Value = 1191  # K93: was 1169.1111
$ws.Cells.Item(93, 13).Value = 57  # M93: was 78.88889999999992
$ws.Cells.Item(122, 8).Value = 3858.7778  # H122: was 4051.8
$ws.Cells.Item(122, 9).Value = 3091.125  # I122: was 3240.682
$ws.Cells.Item(122, 11).Value = 9273.375  # K122: was 9722.045999999998
$ws.Cells.Item(122, 13).Value = -6823.375  # M122: was -7272.045999999998
$ws.Cells.Item(126, 8).Value = 13946.8  # H126: was 12747.429
$ws.Cells.Item(126, 9).Value = 0  # I126: was 9749
$ws.Cells.Item(126, 11).Value = 0  # K126: was 29247
$ws.Cells.Item(126, 13).ClearContents()  # M126: was -26777
$ws.Cells.Item(132, 8).Value = 3389.5686  # H132: was 3299.151
$ws.Cells.Item(132, 9).Value = 3416.4517  # I132: was 3254.4849
$ws.Cells.Item(132, 10).Value = 3347.9  # J132: was 3372.85
$ws.Cells.Item(132, 11).Value = 10249.3551  # K132: was 9763.4547
$ws.Cells.Item(132, 12).Value = 10043.7  # L132: was 10118.55
$ws.Cells.Item(132, 13).Value = -7719.355100000001  # M132: was -7233.4547
$ws.Cells.Item(132, 14).Value = -15103.7  # N132: was -15178.55

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1117.875  # H100: was 1133.4375
$ws.Cells.Item(100, 9).Value = 914.0769  # I100: was 933.2308
$ws.Cells.Item(100, 11).Value = 1828.1538  # K100: was 1866.4616
$ws.Cells.Item(100, 13).Value = -1287.1538  # M100: was -1325.4616
$ws.Cells.Item(107, 8).Value = 1785.0714  # H107: was 2087.4167
$ws.Cells.Item(107, 9).Value = 1586  # I107: was 1879.6666
$ws.Cells.Item(107, 10).Value = 2282.75  # J107: was 2710.6667
$ws.Cells.Item(107, 11).Value = 4758  # K107: was 5638.9998
$ws.Cells.Item(107, 12).Value = 6848.25  # L107: was 8132.000100000001
$ws.Cells.Item(107, 13).Value = -2838  # M107: was -3718.9998
$ws.Cells.Item(107, 14).Value = -10688.25  # N107: was -11972.0001
$ws.Cells.Item(136, 8).Value = 12200976  # H136: was 12505983
$ws.Cells.Item(136, 9).Value = 16133312  # I136: was 16671067
$ws.Cells.Item(136, 11).Value = 48399936  # K136: was 50013201
$ws.Cells.Item(136, 13).Value = -48397386  # M136: was -50010651
